$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.515.15'
$ws.Range('E2').Value = '  +0.05%  '

$ws.Range('D3').Value = '3.124.85'
$ws.Range('E3').Value = '  -1.38%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = '''571.00'
$ws.Range('E5').Value = '  -0.23%  '

$ws.Range('D6').Value = '''163.30'
$ws.Range('E6').Value = '  -3.46%  '

$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('D8').Value = '''0.572'
$ws.Range('E8').Value = '  -5.40%  '

$ws.Range('D9').Value = '3.133.40'
$ws.Range('E9').Value = '  -1.63%  '

$ws.Range('D10').Value = '''0.117'
$ws.Range('E10').Value = '  -1.92%  '

$ws.Range('D11').Value = '''6.61'
$ws.Range('E11').Value = '  -3.08%  '

$ws.Range('D12').Value = '''0.380'
$ws.Range('E12').Value = '  -2.69%  '

$ws.Range('D13').Value = '3.665.65'
$ws.Range('E13').Value = '  -1.53%  '

$ws.Range('E14').Value = '  -2.15%  '

$ws.Range('D15').Value = '64.559.21'
$ws.Range('E15').Value = '  +0.01%  '

$ws.Range('D16').Value = '''24.82'
$ws.Range('E16').Value = '  -2.41%  '

$ws.Range('D17').Value = '3.140.06'
$ws.Range('E17').Value = '  -1.61%  '

$ws.Range('D18').Value = '''0.0000155'
$ws.Range('E18').Value = '  -1.58%  '

$ws.Range('D19').Value = '''407.41'
$ws.Range('E19').Value = '  -2.53%  '

$ws.Range('D20').Value = '''5.21'
$ws.Range('E20').Value = '  -2.16%  '

$ws.Range('D21').Value = '''12.43'
$ws.Range('E21').Value = '  -3.70%  '

$ws.Range('D22').Value = '''7.00'
$ws.Range('E22').Value = '  -2.09%  '

$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  -0.10%  '

$ws.Range('D24').Value = '''68.29'
$ws.Range('E24').Value = '  -2.26%  '

$ws.Range('D25').Value = '''0.480'
$ws.Range('E25').Value = '  -3.75%  '

$ws.Range('D26').Value = '''0.193'
$ws.Range('E26').Value = '  -5.36%  '

$ws.Range('D27').Value = '''0.0000102'
$ws.Range('E27').Value = '  -1.79%  '

$ws.Range('D28').Value = '''9.14'
$ws.Range('E28').Value = '  +3.82%  '

$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  +0.10%  '

$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.11%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.80'
$ws.Range('E31').Value = '  -1.64%  '

$ws.Range('D32').Value = '''21.22'
$ws.Range('E32').Value = '  -2.46%  '

$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '''163.73'
$ws.Range('E33').Value = '  +4.76%  '

$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '''4.93'
$ws.Range('E34').Value = '  -2.98%  '

$ws.Range('D35').Value = '''1.13'
$ws.Range('E35').Value = '  +0.34%  '

$ws.Range('D36').Value = '''6.22'
$ws.Range('E36').Value = '  -2.28%  '

$ws.Range('D37').Value = '''1.35'
$ws.Range('E37').Value = '  -1.01%  '

$ws.Range('D38').Value = '''1.66'
$ws.Range('E38').Value = '  -2.81%  '

$ws.Range('D39').Value = '2.601.70'
$ws.Range('E39').Value = '  -3.94%  '

$ws.Range('D40').Value = '''23.65'
$ws.Range('E40').Value = '  -2.33%  '

$ws.Range('D41').Value = '''4.11'
$ws.Range('E41').Value = '  -2.83%  '

$ws.Range('D42').Value = '''38.23'
$ws.Range('E42').Value = '  -2.11%  '

$ws.Range('D43').Value = '''0.690'
$ws.Range('E43').Value = '  -3.78%  '

$ws.Range('D44').Value = '''0.0616'
$ws.Range('E44').Value = '  -0.79%  '

$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '''5.25'
$ws.Range('E45').Value = '  -4.54%  '

$ws.Range('D46').Value = '''0.0254'
$ws.Range('E46').Value = '  -3.96%  '

$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '''21.16'
$ws.Range('E47').Value = '  -1.79%  '

$ws.Range('D48').Value = '''284.88'
$ws.Range('E48').Value = '  -2.15%  '

$ws.Range('D49').Value = '''0.996'
$ws.Range('E49').Value = '  -0.36%  '

$ws.Range('D50').Value = '''0.0971'
$ws.Range('E50').Value = '  -2.09%  '

$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '''10.48'
$ws.Range('E51').Value = '  +0.22%  '
